$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / label cell stays the same text, rewritten as part of this edit
$ws.Range("A1").Value = "HK_R_acc_G"

# All accuracy values collapse to the same constant figure
$ws.Range("A2:A49").Value = 48.026048026048024
